$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.487.48"
$ws.Range("E2").Value = "  -2.01%  "
$ws.Range("D3").Value = "1.583.33"
$ws.Range("E3").Value = "  -3.24%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.96"
$ws.Range("E5").Value = "  -3.08%  "
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("E8").Value = "  -2.19%  "
$ws.Range("E9").Value = "  -1.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.74"
$ws.Range("E10").Value = "  -2.89%  "
$ws.Range("E11").Value = "  -0.84%  "
$ws.Range("D12").Value = "1.808.51"
$ws.Range("E12").Value = "  -2.97%  "
$ws.Range("D13").Value = "1.596.03"
$ws.Range("E13").Value = "  -2.40%  "
$ws.Range("E14").Value = "  -3.56%  "
$ws.Range("E15").Value = "  -3.92%  "
$ws.Range("D16").Value = "25.518.19"
$ws.Range("E16").Value = "  -1.86%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.24"
$ws.Range("E17").Value = "  -2.20%  "
$ws.Range("D18").Value = "0.0₃0709"
$ws.Range("E18").Value = "  -4.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "187.83"
$ws.Range("E20").Value = "  -1.69%  "
$ws.Range("E21").Value = "  -1.94%  "
$ws.Range("E22").Value = "  -3.74%  "
$ws.Range("E23").Value = "  -2.66%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("E25").Value = "  -3.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.93"
$ws.Range("E26").Value = "  -2.01%  "
$ws.Range("E27").Value = "  -5.23%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "14.95"
$ws.Range("E28").Value = "  -1.53%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.49"
$ws.Range("E29").Value = "  -4.46%  "
$ws.Range("E30").Value = "  -4.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0464"
$ws.Range("E31").Value = "  -3.68%  "
$ws.Range("E32").Value = "  -2.53%  "
$ws.Range("E33").Value = "  -4.15%  "
$ws.Range("E34").Value = "  -0.53%  "
$ws.Range("E35").Value = "  -1.53%  "
$ws.Range("D36").Value = "1.089.02"
$ws.Range("E36").Value = "  -4.09%  "
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.34"
$ws.Range("E38").Value = "  -4.04%  "
$ws.Range("E39").Value = "  -2.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.776"
$ws.Range("E40").Value = "  -10.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.493"
$ws.Range("E41").Value = "  -4.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "95.14"
$ws.Range("E42").Value = "  -3.29%  "
$ws.Range("D43").Value = "1.724.49"
$ws.Range("E43").Value = "  -2.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.07"
$ws.Range("E44").Value = "  -2.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.732"
$ws.Range("E45").Value = "  -5.81%  "
$ws.Range("D46").Value = "0.0₆0108"
$ws.Range("E46").Value = "  -7.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "52.96"
$ws.Range("E47").Value = "  -3.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0507"
$ws.Range("E48").Value = "  -3.83%  "
$ws.Range("E49").Value = "  -4.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.409"
$ws.Range("E50").Value = "  -1.57%  "
$ws.Range("E51").Value = "  -0.11%  "
